# Apply the "5 bytes in packet. 2 curves." edit to the ms_to_hr workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Sampling-rate related values change (250 -> 200 samples/sec) ---
$ws.Range("C1").Value = 200

# --- Packet size change (500 -> 400 bytes) ---
$ws.Range("C4").Value = 400

# --- New column header for the second curve (beats per minute axis) ---
$ws.Range("K4").Value = "уд/мин"

# --- New data points describing the two curves ---
$ws.Range("J5").Value = 60
$ws.Range("K5").Value = 200

$ws.Range("J6").Value = 400
$ws.Range("K6").Value = 30

# Move the active selection to reflect where the user was working.
$ws.Range("K7").Select()

$wb.Application.Calculate()
